$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = [double]35.42516366666666
$ws.Cells.Item(2, 8).Value = [double]106.275491
$ws.Cells.Item(2, 9).Value = [double]0.00832770193000585
$ws.Cells.Item(2, 10).Value = [double]0.008327701930005852
$ws.Cells.Item(2, 13).Value = [double]3.795192333333334
$ws.Cells.Item(2, 14).Value = [double]11.385577
$ws.Cells.Item(2, 15).Value = [double]0.01044213755712683
$ws.Cells.Item(2, 16).Value = [double]0.01044213755712683
$ws.Cells.Item(2, 17).Value = [double]134.4453095548119
$ws.Cells.Item(2, 18).Value = [double]1210.007785993307
$ws.Cells.Item(2, 19).Value = [double]0.0000869590090878717
$ws.Cells.Item(2, 20).Value = [double]0.00008695900908787172
$ws.Cells.Item(3, 7).Value = [double]35.42516366666666
$ws.Cells.Item(3, 8).Value = [double]106.275491
$ws.Cells.Item(3, 9).Value = [double]0.00832770193000585
$ws.Cells.Item(3, 10).Value = [double]0.008327701930005852
$ws.Cells.Item(3, 14).Value = [double]730.1291960000001
$ws.Cells.Item(3, 15).Value = [double]0.6696287328350964
$ws.Cells.Item(3, 16).Value = [double]0.6696287328350964
$ws.Cells.Item(3, 17).Value = [double]8621.648755370581
$ws.Cells.Item(3, 18).Value = [double]77594.83879833524
$ws.Cells.Item(3, 19).Value = [double]0.005576468490818204
$ws.Cells.Item(3, 20).Value = [double]0.005576468490818205
$ws.Cells.Item(4, 7).Value = [double]35.42516366666666
$ws.Cells.Item(4, 8).Value = [double]106.275491
$ws.Cells.Item(4, 9).Value = [double]0.00832770193000585
$ws.Cells.Item(4, 10).Value = [double]0.008327701930005852
$ws.Cells.Item(4, 13).Value = [double]29.801371
$ws.Cells.Item(4, 14).Value = [double]89.404113
$ws.Cells.Item(4, 15).Value = [double]0.08199584844219236
$ws.Cells.Item(4, 16).Value = [double]0.08199584844219235
$ws.Cells.Item(4, 17).Value = [double]1055.718445166053
$ws.Cells.Item(4, 18).Value = [double]9501.466006494482
$ws.Cells.Item(4, 19).Value = [double]0.0006828369853245126
$ws.Cells.Item(4, 20).Value = [double]0.0006828369853245126
$ws.Cells.Item(5, 7).Value = [double]35.42516366666666
$ws.Cells.Item(5, 8).Value = [double]106.275491
$ws.Cells.Item(5, 9).Value = [double]0.00832770193000585
$ws.Cells.Item(5, 10).Value = [double]0.008327701930005852
$ws.Cells.Item(5, 13).Value = [double]86.47679266666667
$ws.Cells.Item(5, 14).Value = [double]259.430378
$ws.Cells.Item(5, 15).Value = [double]0.2379332811655844
$ws.Cells.Item(5, 16).Value = [double]0.2379332811655844
$ws.Cells.Item(5, 17).Value = [double]3063.454533585066
$ws.Cells.Item(5, 18).Value = [double]27571.0908022656
$ws.Cells.Item(5, 19).Value = [double]0.001981437444775262
$ws.Cells.Item(5, 20).Value = [double]0.001981437444775262
$ws.Cells.Item(6, 9).Value = [double]0.01070182047907406
$ws.Cells.Item(6, 10).Value = [double]0.01070182047907406
$ws.Cells.Item(6, 13).Value = [double]3.795192333333334
$ws.Cells.Item(6, 14).Value = [double]11.385577
$ws.Cells.Item(6, 15).Value = [double]0.01044213755712683
$ws.Cells.Item(6, 16).Value = [double]0.01044213755712683
$ws.Cells.Item(6, 17).Value = [double]172.7739031971004
$ws.Cells.Item(6, 18).Value = [double]1554.965128773903
$ws.Cells.Item(6, 19).Value = [double]0.0001117498815541683
$ws.Cells.Item(6, 20).Value = [double]0.0001117498815541683
$ws.Cells.Item(7, 9).Value = [double]0.01070182047907406
$ws.Cells.Item(7, 10).Value = [double]0.01070182047907406
$ws.Cells.Item(7, 14).Value = [double]730.1291960000001
$ws.Cells.Item(7, 15).Value = [double]0.6696287328350964
$ws.Cells.Item(7, 16).Value = [double]0.6696287328350964
$ws.Cells.Item(7, 17).Value = [double]11079.56768735399
$ws.Cells.Item(7, 18).Value = [double]99716.10918618586
$ws.Cells.Item(7, 19).Value = [double]0.007166246486431044
$ws.Cells.Item(7, 20).Value = [double]0.007166246486431044
$ws.Cells.Item(8, 9).Value = [double]0.01070182047907406
$ws.Cells.Item(8, 10).Value = [double]0.01070182047907406
$ws.Cells.Item(8, 13).Value = [double]29.801371
$ws.Cells.Item(8, 14).Value = [double]89.404113
$ws.Cells.Item(8, 15).Value = [double]0.08199584844219236
$ws.Cells.Item(8, 16).Value = [double]0.08199584844219235
$ws.Cells.Item(8, 17).Value = [double]1356.689921370223
$ws.Cells.Item(8, 18).Value = [double]12210.20929233201
$ws.Cells.Item(8, 19).Value = [double]0.0008775048500577067
$ws.Cells.Item(8, 20).Value = [double]0.0008775048500577066
$ws.Cells.Item(9, 9).Value = [double]0.01070182047907406
$ws.Cells.Item(9, 10).Value = [double]0.01070182047907406
$ws.Cells.Item(9, 13).Value = [double]86.47679266666667
$ws.Cells.Item(9, 14).Value = [double]259.430378
$ws.Cells.Item(9, 15).Value = [double]0.2379332811655844
$ws.Cells.Item(9, 16).Value = [double]0.2379332811655844
$ws.Cells.Item(9, 17).Value = [double]3936.805224272705
$ws.Cells.Item(9, 18).Value = [double]35431.24701845434
$ws.Cells.Item(9, 19).Value = [double]0.002546319261031136
$ws.Cells.Item(9, 20).Value = [double]0.002546319261031136
$ws.Cells.Item(10, 7).Value = [double]51.06824600000001
$ws.Cells.Item(10, 8).Value = [double]153.204738
$ws.Cells.Item(10, 9).Value = [double]0.01200505761322374
$ws.Cells.Item(10, 10).Value = [double]0.01200505761322374
$ws.Cells.Item(10, 13).Value = [double]3.795192333333334
$ws.Cells.Item(10, 14).Value = [double]11.385577
$ws.Cells.Item(10, 15).Value = [double]0.01044213755712683
$ws.Cells.Item(10, 16).Value = [double]0.01044213755712683
$ws.Cells.Item(10, 17).Value = [double]193.8138156959807
$ws.Cells.Item(10, 18).Value = [double]1744.324341263826
$ws.Cells.Item(10, 19).Value = [double]0.000125358462978515
$ws.Cells.Item(10, 20).Value = [double]0.000125358462978515
$ws.Cells.Item(11, 7).Value = [double]51.06824600000001
$ws.Cells.Item(11, 8).Value = [double]153.204738
$ws.Cells.Item(11, 9).Value = [double]0.01200505761322374
$ws.Cells.Item(11, 10).Value = [double]0.01200505761322374
$ws.Cells.Item(11, 14).Value = [double]730.1291960000001
$ws.Cells.Item(11, 15).Value = [double]0.6696287328350964
$ws.Cells.Item(11, 16).Value = [double]0.6696287328350964
$ws.Cells.Item(11, 17).Value = [double]12428.80579770341
$ws.Cells.Item(11, 18).Value = [double]111859.2521793307
$ws.Cells.Item(11, 19).Value = [double]0.008038931517155339
$ws.Cells.Item(11, 20).Value = [double]0.00803893151715534
$ws.Cells.Item(12, 7).Value = [double]51.06824600000001
$ws.Cells.Item(12, 8).Value = [double]153.204738
$ws.Cells.Item(12, 9).Value = [double]0.01200505761322374
$ws.Cells.Item(12, 10).Value = [double]0.01200505761322374
$ws.Cells.Item(12, 13).Value = [double]29.801371
$ws.Cells.Item(12, 14).Value = [double]89.404113
$ws.Cells.Item(12, 15).Value = [double]0.08199584844219236
$ws.Cells.Item(12, 16).Value = [double]0.08199584844219235
$ws.Cells.Item(12, 17).Value = [double]1521.903745365266
$ws.Cells.Item(12, 18).Value = [double]13697.13370828739
$ws.Cells.Item(12, 19).Value = [double]0.0009843648845936812
$ws.Cells.Item(12, 20).Value = [double]0.0009843648845936812
$ws.Cells.Item(13, 7).Value = [double]51.06824600000001
$ws.Cells.Item(13, 8).Value = [double]153.204738
$ws.Cells.Item(13, 9).Value = [double]0.01200505761322374
$ws.Cells.Item(13, 10).Value = [double]0.01200505761322374
$ws.Cells.Item(13, 13).Value = [double]86.47679266666667
$ws.Cells.Item(13, 14).Value = [double]259.430378
$ws.Cells.Item(13, 15).Value = [double]0.2379332811655844
$ws.Cells.Item(13, 16).Value = [double]0.2379332811655844
$ws.Cells.Item(13, 17).Value = [double]4416.21812119233
$ws.Cells.Item(13, 18).Value = [double]39745.96309073097
$ws.Cells.Item(13, 19).Value = [double]0.002856402748496203
$ws.Cells.Item(13, 20).Value = [double]0.002856402748496203
$ws.Cells.Item(14, 7).Value = [double]4121.876464666667
$ws.Cells.Item(14, 8).Value = [double]12365.629394
$ws.Cells.Item(14, 9).Value = [double]0.9689654199776964
$ws.Cells.Item(14, 10).Value = [double]0.9689654199776964
$ws.Cells.Item(14, 13).Value = [double]3.795192333333334
$ws.Cells.Item(14, 14).Value = [double]11.385577
$ws.Cells.Item(14, 15).Value = [double]0.01044213755712683
$ws.Cells.Item(14, 16).Value = [double]0.01044213755712683
$ws.Cells.Item(14, 17).Value = [double]15643.31395765004
$ws.Cells.Item(14, 18).Value = [double]140789.8256188504
$ws.Cells.Item(14, 19).Value = [double]0.01011807020350628
$ws.Cells.Item(14, 20).Value = [double]0.01011807020350628
$ws.Cells.Item(15, 7).Value = [double]4121.876464666667
$ws.Cells.Item(15, 8).Value = [double]12365.629394
$ws.Cells.Item(15, 9).Value = [double]0.9689654199776964
$ws.Cells.Item(15, 10).Value = [double]0.9689654199776964
$ws.Cells.Item(15, 14).Value = [double]730.1291960000001
$ws.Cells.Item(15, 15).Value = [double]0.6696287328350964
$ws.Cells.Item(15, 16).Value = [double]0.6696287328350964
$ws.Cells.Item(15, 17).Value = [double]1003167.449719465
$ws.Cells.Item(15, 18).Value = [double]9028507.047475189
$ws.Cells.Item(15, 19).Value = [double]0.6488470863406918
$ws.Cells.Item(15, 20).Value = [double]0.6488470863406918
$ws.Cells.Item(16, 7).Value = [double]4121.876464666667
$ws.Cells.Item(16, 8).Value = [double]12365.629394
$ws.Cells.Item(16, 9).Value = [double]0.9689654199776964
$ws.Cells.Item(16, 10).Value = [double]0.9689654199776964
$ws.Cells.Item(16, 13).Value = [double]29.801371
$ws.Cells.Item(16, 14).Value = [double]89.404113
$ws.Cells.Item(16, 15).Value = [double]0.08199584844219236
$ws.Cells.Item(16, 16).Value = [double]0.08199584844219235
$ws.Cells.Item(16, 17).Value = [double]122837.5697396997
$ws.Cells.Item(16, 18).Value = [double]1105538.127657298
$ws.Cells.Item(16, 19).Value = [double]0.07945114172221646
$ws.Cells.Item(16, 20).Value = [double]0.07945114172221644
$ws.Cells.Item(17, 7).Value = [double]4121.876464666667
$ws.Cells.Item(17, 8).Value = [double]12365.629394
$ws.Cells.Item(17, 9).Value = [double]0.9689654199776964
$ws.Cells.Item(17, 10).Value = [double]0.9689654199776964
$ws.Cells.Item(17, 13).Value = [double]86.47679266666667
$ws.Cells.Item(17, 14).Value = [double]259.430378
$ws.Cells.Item(17, 15).Value = [double]0.2379332811655844
$ws.Cells.Item(17, 16).Value = [double]0.2379332811655844
$ws.Cells.Item(17, 17).Value = [double]356446.6564325924
$ws.Cells.Item(17, 18).Value = [double]3208019.907893331
$ws.Cells.Item(17, 19).Value = [double]0.2305491217112818
$ws.Cells.Item(17, 20).Value = [double]0.2305491217112818
